# Apply updated crypto prices / 1h volume changes (and a few reordered rows)
# as described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.967.77"
$ws.Range("E2").Value = "'  +0.48%  "
$ws.Range("D3").Value = "'2.262.67"
$ws.Range("E3").Value = "'  -0.56%  "
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("B5").Value = "'XRP"
$ws.Range("C5").Value = "'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "'0.652"
$ws.Range("E5").Value = "'  +4.20%  "
$ws.Range("B6").Value = "'BNB"
$ws.Range("C6").Value = "'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'232.93"
$ws.Range("E6").Value = "'  +0.72%  "
$ws.Range("D7").Value = "'63.55"
$ws.Range("E7").Value = "'  -0.44%  "
$ws.Range("E8").Value = "'  -0.10%  "
$ws.Range("E9").Value = "'  +4.85%  "
$ws.Range("D10").Value = "'0.0972"
$ws.Range("E10").Value = "'  -7.15%  "
$ws.Range("D11").Value = "'58.09"
$ws.Range("E11").Value = "'  +1.30%  "
$ws.Range("E12").Value = "'  +2.45%  "
$ws.Range("E13").Value = "'  +1.39%  "
$ws.Range("D14").Value = "'2.597.28"
$ws.Range("E14").Value = "'  -0.49%  "
$ws.Range("E15").Value = "'  -0.29%  "
$ws.Range("D16").Value = "'6.15"
$ws.Range("E16").Value = "'  +4.66%  "
$ws.Range("D17").Value = "'0.843"
$ws.Range("E17").Value = "'  +3.04%  "
$ws.Range("D18").Value = "'2.269.41"
$ws.Range("E18").Value = "'  -0.29%  "
$ws.Range("D19").Value = "'43.824.55"
$ws.Range("E19").Value = "'  +0.43%  "
$ws.Range("D20").Value = "'0.0₃0979"
$ws.Range("E20").Value = "'  -3.79%  "
$ws.Range("D21").Value = "'73.74"
$ws.Range("E21").Value = "'  +0.55%  "
$ws.Range("D22").Value = "'6.17"
$ws.Range("E22").Value = "'  +1.28%  "
$ws.Range("D23").Value = "'249.54"
$ws.Range("E23").Value = "'  +0.23%  "
$ws.Range("E24").Value = "'  -0.11%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "'  -1.59%  "
$ws.Range("B26").Value = "'Toncoin"
$ws.Range("C26").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "'  +3.92%  "
$ws.Range("B27").Value = "'WEMIXToken"
$ws.Range("C27").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "'3.52"
$ws.Range("E27").Value = "'  +26.22%  "
$ws.Range("D28").Value = "'9.91"
$ws.Range("E28").Value = "'  +0.88%  "
$ws.Range("D29").Value = "'173.85"
$ws.Range("E29").Value = "'  +1.18%  "
$ws.Range("D30").Value = "'21.92"
$ws.Range("E30").Value = "'  +3.96%  "
$ws.Range("E31").Value = "'  -0.38%  "
$ws.Range("D32").Value = "'1.43"
$ws.Range("E32").Value = "'  -0.73%  "
$ws.Range("E33").Value = "'  +3.67%  "
$ws.Range("D34").Value = "'4.96"
$ws.Range("E34").Value = "'  +5.63%  "
$ws.Range("D35").Value = "'0.0687"
$ws.Range("E35").Value = "'  -0.31%  "
$ws.Range("D36").Value = "'4.95"
$ws.Range("E36").Value = "'  -2.57%  "
$ws.Range("E37").Value = "'  -2.94%  "
$ws.Range("E38").Value = "'  -4.96%  "
$ws.Range("E39").Value = "'  -1.61%  "
$ws.Range("D40").Value = "'0.0256"
$ws.Range("E40").Value = "'  +3.24%  "
$ws.Range("E41").Value = "'  +0.13%  "
$ws.Range("D42").Value = "'8.68"
$ws.Range("E42").Value = "'  +3.69%  "
$ws.Range("D43").Value = "'17.26"
$ws.Range("E43").Value = "'  +0.71%  "
$ws.Range("D44").Value = "'98.77"
$ws.Range("E44").Value = "'  +1.45%  "
$ws.Range("B45").Value = "'Cronos"
$ws.Range("C45").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0950"
$ws.Range("E45").Value = "'  -1.11%  "
$ws.Range("B46").Value = "'TrustWalletToken"
$ws.Range("C46").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.19"
$ws.Range("E46").Value = "'  -0.96%  "
$ws.Range("B47").Value = "'FTXToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "'4.38"
$ws.Range("E47").Value = "'  -0.18%  "
$ws.Range("B48").Value = "'TerraClassic"
$ws.Range("C48").Value = "'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D48").Value = "'0.000209"
$ws.Range("E48").Value = "'  +2.88%  "
$ws.Range("D49").Value = "'1.456.49"
$ws.Range("E49").Value = "'  -1.34%  "
$ws.Range("D50").Value = "'2.34"
$ws.Range("E50").Value = "'  +0.61%  "
$ws.Range("D51").Value = "'9.96"
$ws.Range("E51").Value = "'  -5.10%  "
